$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp label (A1)
$ws.Range("A1").Value = "Datos actualizados a 5 de Julio de 2020 a las 02:37"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 2935427
$ws.Range("C4").Value = 44839
$ws.Range("D4").Value = 1259468
$ws.Range("E4").Value = 1543646
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 249
$ws.Range("H4").Value = 132313

# Row 23: Canada
$ws.Range("A23").Value = "Canada"
$ws.Range("B23").Value = 105317
$ws.Range("C23").Value = 226
$ws.Range("D23").Value = 68990
$ws.Range("E23").Value = 27653
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 11
$ws.Range("H23").Value = 8674

# Row 26: Argentina
$ws.Range("A26").Value = "Argentina"
$ws.Range("B26").Value = 75376
$ws.Range("C26").Value = 2590
$ws.Range("D26").Value = 25930
$ws.Range("E26").Value = 47965
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 44
$ws.Range("H26").Value = 1481

# Row 27: Egipto
$ws.Range("A27").Value = "Egipto"
$ws.Range("B27").Value = 74035
$ws.Range("C27").Value = 1324
$ws.Range("D27").Value = 20103
$ws.Range("E27").Value = 50652
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 79
$ws.Range("H27").Value = 3280

# Row 43: Panama
$ws.Range("A43").Value = "Panama"
$ws.Range("B43").Value = 36983
$ws.Range("C43").Value = 988
$ws.Range("D43").Value = 17761
$ws.Range("E43").Value = 18502
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 22
$ws.Range("H43").Value = 720

# Row 44: Bolivia
$ws.Range("A44").Value = "Bolivia"
$ws.Range("B44").Value = 36818
$ws.Range("C44").Value = 1290
$ws.Range("D44").Value = 10766
$ws.Range("E44").Value = 24732
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 49
$ws.Range("H44").Value = 1320

# Row 45: Republica Dominicana
$ws.Range("A45").Value = "Republica Dominicana"
$ws.Range("B45").Value = 36184
$ws.Range("C45").Value = 1036
$ws.Range("D45").Value = 18602
$ws.Range("E45").Value = 16796
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 11
$ws.Range("H45").Value = 786

# Row 122: Malaui
$ws.Range("A122").Value = "Malaui"
$ws.Range("B122").Value = 1613
$ws.Range("C122").Value = 115
$ws.Range("D122").Value = 317
$ws.Range("E122").Value = 1279
$ws.Range("F122").Value = 0
$ws.Range("G122").Value = 1
$ws.Range("H122").Value = 17

# Row 123: Congo
$ws.Range("A123").Value = "Congo"
$ws.Range("B123").Value = 1557
$ws.Range("C123").Value = 0
$ws.Range("D123").Value = 501
$ws.Range("E123").Value = 1012
$ws.Range("F123").Value = 0
$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 44

# Row 124: Sierra Leona
$ws.Range("A124").Value = "Sierra Leona"
$ws.Range("B124").Value = 1533
$ws.Range("C124").Value = 9
$ws.Range("D124").Value = 1051
$ws.Range("E124").Value = 420
$ws.Range("F124").Value = 0
$ws.Range("G124").Value = 0
$ws.Range("H124").Value = 62

# Row 125: Nueva Zelanda
$ws.Range("A125").Value = "Nueva Zelanda"
$ws.Range("B125").Value = 1530
$ws.Range("C125").Value = 0
$ws.Range("D125").Value = 1490
$ws.Range("E125").Value = 18
$ws.Range("F125").Value = 0
$ws.Range("G125").Value = 0
$ws.Range("H125").Value = 22

# Row 136: Libia
$ws.Range("A136").Value = "Libia"
$ws.Range("B136").Value = 989
$ws.Range("C136").Value = 71
$ws.Range("D136").Value = 258
$ws.Range("E136").Value = 704
$ws.Range("F136").Value = 0
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 27

# Row 137: Burkina Faso
$ws.Range("A137").Value = "Burkina Faso"
$ws.Range("B137").Value = 987
$ws.Range("C137").Value = 7
$ws.Range("D137").Value = 854
$ws.Range("E137").Value = 80
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 53

# Row 138: Mozambique
$ws.Range("A138").Value = "Mozambique"
$ws.Range("B138").Value = 969
$ws.Range("C138").Value = 30
$ws.Range("D138").Value = 256
$ws.Range("E138").Value = 706
$ws.Range("F138").Value = 0
$ws.Range("G138").Value = 1
$ws.Range("H138").Value = 7

# Row 139: Uruguay
$ws.Range("A139").Value = "Uruguay"
$ws.Range("B139").Value = 955
$ws.Range("C139").Value = 3
$ws.Range("D139").Value = 840
$ws.Range("E139").Value = 87
$ws.Range("F139").Value = 0
$ws.Range("G139").Value = 0
$ws.Range("H139").Value = 28

# Row 140: Suazilandia
$ws.Range("A140").Value = "Suazilandia"
$ws.Range("B140").Value = 954
$ws.Range("C140").Value = 45
$ws.Range("D140").Value = 535
$ws.Range("E140").Value = 406
$ws.Range("F140").Value = 0
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 13

# Row 141: Georgia
$ws.Range("A141").Value = "Georgia"
$ws.Range("B141").Value = 948
$ws.Range("C141").Value = 5
$ws.Range("D141").Value = 825
$ws.Range("E141").Value = 108
$ws.Range("F141").Value = 0
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = 15

# Row 142: Uganda
$ws.Range("A142").Value = "Uganda"
$ws.Range("B142").Value = 927
$ws.Range("C142").Value = 16
$ws.Range("D142").Value = 868
$ws.Range("E142").Value = 59
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 0

# Row 152: Togo
$ws.Range("A152").Value = "Togo"
$ws.Range("B152").Value = 676
$ws.Range("C152").Value = 5
$ws.Range("D152").Value = 432
$ws.Range("E152").Value = 229
$ws.Range("F152").Value = 0
$ws.Range("G152").Value = 1
$ws.Range("H152").Value = 15

# Row 153: Malta
$ws.Range("A153").Value = "Malta"
$ws.Range("B153").Value = 672
$ws.Range("C153").Value = 0
$ws.Range("D153").Value = 651
$ws.Range("E153").Value = 12
$ws.Range("F153").Value = 0
$ws.Range("G153").Value = 0
$ws.Range("H153").Value = 9

# Row 155: Reunion
$ws.Range("A155").Value = "Reunion"
$ws.Range("B155").Value = 536
$ws.Range("C155").Value = 3
$ws.Range("D155").Value = 472
$ws.Range("E155").Value = 62
$ws.Range("F155").Value = 0
$ws.Range("G155").Value = 0
$ws.Range("H155").Value = 2

# Row 167: Guyana
$ws.Range("A167").Value = "Guyana"
$ws.Range("B167").Value = 272
$ws.Range("C167").Value = 16
$ws.Range("D167").Value = 120
$ws.Range("E167").Value = 138
$ws.Range("F167").Value = 0
$ws.Range("G167").Value = 0
$ws.Range("H167").Value = 14

# Row 205: Fiyi
$ws.Range("A205").Value = "Fiyi"
$ws.Range("B205").Value = 18
$ws.Range("C205").Value = 0
$ws.Range("D205").Value = 18
$ws.Range("E205").Value = 0
$ws.Range("F205").Value = 0
$ws.Range("G205").Value = 0
$ws.Range("H205").Value = 0

# Row 206: Dominica
$ws.Range("A206").Value = "Dominica"
$ws.Range("B206").Value = 18
$ws.Range("C206").Value = 0
$ws.Range("D206").Value = 18
$ws.Range("E206").Value = 0
$ws.Range("F206").Value = 0
$ws.Range("G206").Value = 0
$ws.Range("H206").Value = 0
